$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# row 43 (hunk 0)
$ws.Range("H43").Value = 1131.1818
$ws.Range("I43").Value = 696.5
$ws.Range("J43").Value = 1227.7778
$ws.Range("K43").Value = 696.5
$ws.Range("L43").Value = 1227.7778
$ws.Range("M43").Value = -627.5
$ws.Range("N43").Value = -1365.7778
# row 64 (hunk 1)
$ws.Range("H64").Value = 3466.6667
$ws.Range("J64").Value = 3600
$ws.Range("L64").Value = 3600
$ws.Range("N64").Value = -4096
# row 67 (hunk 2)
$ws.Range("H67").Value = 3466.6667
$ws.Range("J67").Value = 3600
$ws.Range("L67").Value = 3600
$ws.Range("N67").Value = -5316
# row 74 (hunk 3)
$ws.Range("H74").Value = 4618.5
$ws.Range("I74").Value = 3350
$ws.Range("K74").Value = 3350
$ws.Range("M74").Value = -2414
# row 76 (hunk 4)
$ws.Range("H76").Value = 3200
$ws.Range("I76").Value = 3212.9033
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 3212.9033
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -2897.9033
$ws.Range("N76").Value = -3630
# row 77 (hunk 5)
$ws.Range("H77").Value = 4618.5
$ws.Range("I77").Value = 3350
$ws.Range("K77").Value = 16750
$ws.Range("M77").Value = -12070
# row 79 (hunk 6)
$ws.Range("H79").Value = 3200
$ws.Range("I79").Value = 3212.9033
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 3212.9033
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -2120.9033
$ws.Range("N79").Value = -5184
# row 98 (hunk 7)
$ws.Range("H98").Value = 1017.5
$ws.Range("I98").Value = 276.15384
$ws.Range("J98").Value = 2088.3333
$ws.Range("K98").Value = 276.15384
$ws.Range("L98").Value = 2088.3333
$ws.Range("M98").Value = 1221.84616
$ws.Range("N98").Value = -5084.3333
# row 112 (hunk 8)
$ws.Range("H112").Value = 4417.0557
$ws.Range("I112").Value = 1550
$ws.Range("K112").Value = 4650
$ws.Range("M112").Value = -3542
# row 116 (hunk 9)
$ws.Range("H116").Value = 3097.3333
$ws.Range("I116").Value = 2907
$ws.Range("K116").Value = 2907
$ws.Range("M116").Value = 535
# row 122 (hunk 10)
$ws.Range("H122").Value = 1017.5
$ws.Range("I122").Value = 276.15384
$ws.Range("J122").Value = 2088.3333
$ws.Range("K122").Value = 828.4615200000001
$ws.Range("L122").Value = 6264.999899999999
$ws.Range("M122").Value = 1621.53848
$ws.Range("N122").Value = -11164.9999
# row 132 (hunk 11)
$ws.Range("H132").Value = 5548.037
$ws.Range("I132").Value = 5371.3335
$ws.Range("J132").Value = 6166.5
$ws.Range("K132").Value = 16114.0005
$ws.Range("L132").Value = 18499.5
$ws.Range("M132").Value = -13584.0005
$ws.Range("N132").Value = -23559.5
# row 135 (hunk 12)
$ws.Range("H135").Value = 4785.1763
$ws.Range("I135").Value = 3592
$ws.Range("J135").Value = 6972.6665
$ws.Range("K135").Value = 32328
$ws.Range("L135").Value = 62753.9985
$ws.Range("M135").Value = -29793
$ws.Range("N135").Value = -67823.9985
# row 138 (hunk 13)
$ws.Range("H138").Value = 2317.772
$ws.Range("I138").Value = 3122.7693
$ws.Range("J138").Value = 2079.932
$ws.Range("K138").Value = 9368.3079
$ws.Range("L138").Value = 6239.795999999999
$ws.Range("M138").Value = -4228.3079
$ws.Range("N138").Value = -16519.796
# row 141 (hunk 14)
$ws.Range("H141").Value = 7400.364
$ws.Range("I141").Value = 2663.7693
$ws.Range("J141").Value = 14242.111
$ws.Range("K141").Value = 7991.3079
$ws.Range("L141").Value = 42726.333
$ws.Range("M141").Value = -2811.3079
$ws.Range("N141").Value = -53086.333

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# row 32 (hunk 15)
$ws.Range("H32").Value = 419435.22
$ws.Range("I32").Value = 471139.75
$ws.Range("K32").Value = 471139.75
$ws.Range("M32").Value = -470852.75
# row 63 (hunk 16)
$ws.Range("H63").Value = 5887.7
$ws.Range("I63").Value = 2666.6667
$ws.Range("J63").Value = 7268.143
$ws.Range("K63").Value = 2666.6667
$ws.Range("L63").Value = 7268.143
$ws.Range("M63").Value = -1980.6667
$ws.Range("N63").Value = -8640.143
# row 66 (hunk 17)
$ws.Range("H66").Value = 5887.7
$ws.Range("I66").Value = 2666.6667
$ws.Range("J66").Value = 7268.143
$ws.Range("K66").Value = 13333.3335
$ws.Range("L66").Value = 36340.715
$ws.Range("M66").Value = -9901.333500000001
$ws.Range("N66").Value = -43204.715
# row 74 (hunk 18)
$ws.Range("H74").Value = 1972.7858
$ws.Range("I74").Value = 1454
$ws.Range("J74").Value = 2735.7058
$ws.Range("K74").Value = 1454
$ws.Range("L74").Value = 2735.7058
$ws.Range("M74").Value = -580
$ws.Range("N74").Value = -4483.7058
# row 77 (hunk 19)
$ws.Range("H77").Value = 1972.7858
$ws.Range("I77").Value = 1454
$ws.Range("J77").Value = 2735.7058
$ws.Range("K77").Value = 7270
$ws.Range("L77").Value = 13678.529
$ws.Range("M77").Value = -2902
$ws.Range("N77").Value = -22414.529
# row 122 (hunk 20)
$ws.Range("H122").Value = 61399.883
$ws.Range("I122").Value = 112622.11
$ws.Range("K122").Value = 337866.33
$ws.Range("M122").Value = -335416.33

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# row 31 (hunk 21)
$ws.Range("H31").Value = 6127.213
$ws.Range("I31").Value = 1573.125
$ws.Range("J31").Value = 10879.305
$ws.Range("K31").Value = 1573.125
$ws.Range("L31").Value = 10879.305
$ws.Range("M31").Value = -1278.125
$ws.Range("N31").Value = -11469.305
# row 34 (hunk 22)
$ws.Range("H34").Value = 6127.213
$ws.Range("I34").Value = 1573.125
$ws.Range("J34").Value = 10879.305
$ws.Range("K34").Value = 1573.125
$ws.Range("L34").Value = 10879.305
$ws.Range("M34").Value = -1371.125
$ws.Range("N34").Value = -11283.305
# row 58 (hunk 23)
$ws.Range("H58").Value = 993.775
$ws.Range("I58").Value = 692.11536
$ws.Range("K58").Value = 692.11536
$ws.Range("M58").Value = -489.11536
# row 62 (hunk 24)
$ws.Range("H62").Value = 4666.6665
$ws.Range("I62").Value = 4600
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 4600
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -3976
$ws.Range("N62").Value = -6248
# row 65 (hunk 25)
$ws.Range("H65").Value = 4666.6665
$ws.Range("I65").Value = 4600
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 23000
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -19880
$ws.Range("N65").Value = -31240
# row 132 (hunk 26)
$ws.Range("H132").Value = 9805714
$ws.Range("I132").Value = 1168.3636
$ws.Range("J132").Value = 27780714
$ws.Range("K132").Value = 3505.0908
$ws.Range("L132").Value = 83342142
$ws.Range("M132").Value = -975.0907999999999
$ws.Range("N132").Value = -83347202
# row 134 (hunk 27)
$ws.Range("H134").Value = 2486.6667
$ws.Range("I134").Value = 1572.8948
$ws.Range("J134").Value = 5959
$ws.Range("K134").Value = 4718.6844
$ws.Range("L134").Value = 17877
$ws.Range("M134").Value = -2183.6844
$ws.Range("N134").Value = -22947
# row 136 (hunk 28)
$ws.Range("H136").Value = 993.775
$ws.Range("I136").Value = 692.11536
$ws.Range("K136").Value = 2076.34608
$ws.Range("M136").Value = 473.6539199999997

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# row 5 (hunk 29)
$ws.Range("H5").Value = 539.65
$ws.Range("I5").Value = 458.41177
$ws.Range("K5").Value = 1375.23531
$ws.Range("M5").Value = -1263.23531
# row 131 (hunk 30)
$ws.Range("H131").Value = 980.8095
$ws.Range("I131").Value = 412.5
$ws.Range("J131").Value = 1114.5294
$ws.Range("K131").Value = 1237.5
$ws.Range("L131").Value = 3343.5882
$ws.Range("M131").Value = 3802.5
$ws.Range("N131").Value = -13423.5882
# row 132 (hunk 31)
$ws.Range("H132").Value = 2526.6875
$ws.Range("I132").Value = 2309
$ws.Range("K132").Value = 20781
$ws.Range("M132").Value = -18251
# row 135 (hunk 32)
$ws.Range("H135").Value = 539.65
$ws.Range("I135").Value = 458.41177
$ws.Range("K135").Value = 4125.70593
$ws.Range("M135").Value = -1590.70593

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# row 70 (hunk 33)
$ws.Range("H70").Value = 5662.9473
$ws.Range("I70").Value = 5748
$ws.Range("J70").Value = 5532.533
$ws.Range("K70").Value = 5748
$ws.Range("L70").Value = 5532.533
$ws.Range("M70").Value = -5478
$ws.Range("N70").Value = -6072.533
# row 73 (hunk 34)
$ws.Range("H73").Value = 5662.9473
$ws.Range("I73").Value = 5748
$ws.Range("J73").Value = 5532.533
$ws.Range("K73").Value = 5748
$ws.Range("L73").Value = 5532.533
$ws.Range("M73").Value = -4812
$ws.Range("N73").Value = -7404.533
# row 80 (hunk 35)
$ws.Range("H80").Value = 629980.6
$ws.Range("I80").Value = 1804602
$ws.Range("J80").Value = 42669.9
$ws.Range("K80").Value = 1804602
$ws.Range("L80").Value = 42669.9
$ws.Range("M80").Value = -1803604
$ws.Range("N80").Value = -44665.9
# row 83 (hunk 36)
$ws.Range("H83").Value = 629980.6
$ws.Range("I83").Value = 1804602
$ws.Range("J83").Value = 42669.9
$ws.Range("K83").Value = 9023010
$ws.Range("L83").Value = 213349.5
$ws.Range("M83").Value = -9018018
$ws.Range("N83").Value = -223333.5
# row 126 (hunk 37)
$ws.Range("H126").Value = 2152.2942
$ws.Range("I126").Value = 1925
$ws.Range("J126").Value = 2569
$ws.Range("K126").Value = 5775
$ws.Range("L126").Value = 7707
$ws.Range("M126").Value = -3305
$ws.Range("N126").Value = -12647
# row 132 (hunk 38)
$ws.Range("H132").Value = 2117.6052
$ws.Range("I132").Value = 1652.8276
$ws.Range("K132").Value = 4958.4828
$ws.Range("M132").Value = -2428.4828

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# row 7 (hunk 39)
$ws.Range("H7").Value = 1375.6666
$ws.Range("I7").Value = 1045.7273
$ws.Range("J7").Value = 5005
$ws.Range("K7").Value = 1045.7273
$ws.Range("L7").Value = 5005
$ws.Range("M7").Value = -933.7273
$ws.Range("N7").Value = -5229
# row 16 (hunk 40)
$ws.Range("H16").Value = 901.6667
$ws.Range("I16").Value = 701.9375
$ws.Range("K16").Value = 701.9375
$ws.Range("M16").Value = -531.9375
# row 22 (hunk 41)
$ws.Range("H22").Value = 3909.0571
$ws.Range("I22").Value = 404.7647
$ws.Range("J22").Value = 7218.6665
$ws.Range("K22").Value = 404.7647
$ws.Range("L22").Value = 7218.6665
$ws.Range("M22").Value = -109.7647
$ws.Range("N22").Value = -7808.6665
# row 27 (hunk 42)
$ws.Range("H27").Value = 3909.0571
$ws.Range("I27").Value = 404.7647
$ws.Range("J27").Value = 7218.6665
$ws.Range("K27").Value = 404.7647
$ws.Range("L27").Value = 7218.6665
$ws.Range("M27").Value = -297.7647
$ws.Range("N27").Value = -7432.6665
# row 126 (hunk 43)
$ws.Range("H126").Value = 1375.6666
$ws.Range("I126").Value = 1045.7273
$ws.Range("J126").Value = 5005
$ws.Range("K126").Value = 3137.1819
$ws.Range("L126").Value = 15015
$ws.Range("M126").Value = -667.1819
$ws.Range("N126").Value = -19955

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# row 62 (hunk 44)
$ws.Range("H62").Value = 151250
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 151250
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 151250
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -152498
# row 65 (hunk 45)
$ws.Range("H65").Value = 151250
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 151250
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 756250
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -762490
